# C5-PowerPoint.pptx edit
#  1. Slide 6's table switches to a different (built-in) table style.
#  2. The deck's theme colour scheme is changed from the custom
#     "Integral" palette to the stock Office default palette.

function HexToComRGB($hex) {
    # PowerPoint's RGB COM property stores colours as 0x00BBGGRR, so the
    # byte order needs flipping relative to the "RRGGBB" hex most people
    # think in.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{9DA77A84-6EFB-4549-89F7-42DAD366F3DE}", $false)
    }
}

# --- 2. Theme colours: Integral -> Office default -------------------------
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToComRGB($officeColors[$i - 1])
}
